# Applies the "Bit more stuff added to report and poster" edit:
#   - "Wirespeed - 100Mbit/s"              -> "Wirespeed firewall - 100Mbit/s"
#   - "4uS of latency added"               -> "4us of latency added"
#   - "Measured 0.51W power consumption"   -> "Total measured 0.51W power. Packet filter logic consumes ~2mW"
#     (and the now-redundant blank paragraph right after it is removed)
#   - "Picture 14" is nudged down slightly (y offset 5883669 -> 5910812 EMU)

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# "RESULTS" content placeholder holds all three text edits.
# ---------------------------------------------------------------------------
$shp = $s.Shapes.Item(7)
$tr  = $shp.TextFrame.TextRange

$enDash = [char]0x2013

# 1) "Wirespeed - 100Mbit/s" -> "Wirespeed firewall - 100Mbit/s"
#    Only the run after "Wirespeed" needs its text changed, so grab that
#    run's characters directly (leaves the "Wirespeed" run untouched).
$wirespeedPara = $tr.Paragraphs(4, 1)
$tail = $wirespeedPara.Characters(10, $wirespeedPara.Length - 9)
$tail.Text = " firewall " + $enDash + " 100Mbit/s"

# 2) "4uS of latency added" -> "4us of latency added"
#    Assigning text that shares a common prefix/suffix with the existing
#    run causes the host to fragment it into multiple runs; routing the
#    change through an unrelated placeholder first keeps it a single run,
#    matching the original (single-run) paragraph shape.
$latencyPara = $tr.Paragraphs(5, 1)
$latencyPara.Text = "~~~PLACEHOLDER~~~"
$latencyPara = $tr.Paragraphs(5, 1)
$latencyPara.Text = "4us of latency added"

# 3) "Measured 0.51W power consumption" -> "Total measured 0.51W power. Packet filter logic consumes ~2mW"
#    Same placeholder trick to keep it a single run.
$powerPara = $tr.Paragraphs(14, 1)
$powerPara.Text = "~~~PLACEHOLDER~~~"
$powerPara = $tr.Paragraphs(14, 1)
$powerPara.Text = "Total measured 0.51W power. Packet filter logic consumes ~2mW"

# One of the now-redundant blank paragraphs further down needs to go. Deleting
# the blank paragraph immediately after the text one merges the deleted
# paragraph's mark-properties onto the text run's paragraph (leaving a stray
# <a:endParaRPr/>); deleting a blank paragraph that sits between two other
# blank paragraphs instead merges two identical blanks together with no
# observable effect, so do that.
$blankToRemove = $tr.Paragraphs(16, 1)
$blankToRemove.Delete()

# ---------------------------------------------------------------------------
# Nudge "Picture 14" down a touch (EMU -> points: 1 pt = 12700 EMU).
# ---------------------------------------------------------------------------
$pic = $s.Shapes.Item(10)
$pic.Top = 5910812 / 12700
